$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title text updates (Volume/Number and report week date range) ---
$ws.Range("A8").Value = "Volume 30   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/15/2023  Through  5/21/2023"

# --- Weekly crime statistics table updates (rows 15-30) ---

$ws.Range("D15").Value = 2
$ws.Range("E15").Value = -50
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = -60
$ws.Range("I15").Value = 4
$ws.Range("J15").Value = 10
$ws.Range("K15").Value = -60
$ws.Range("L15").Value = -33.333333333333
$ws.Range("M15").Value = 100
$ws.Range("N15").Value = -50
$ws.Range("C16").Value = 8
$ws.Range("D16").Value = 13
$ws.Range("E16").Value = -38.461538461538
$ws.Range("F16").Value = 36
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 193
$ws.Range("J16").Value = 207
$ws.Range("K16").Value = -6.763285024154
$ws.Range("L16").Value = 39.855072463768
$ws.Range("M16").Value = 238.59649122807
$ws.Range("N16").Value = -80.225409836065
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 14
$ws.Range("E17").Value = -78.571428571428
$ws.Range("F17").Value = 23
$ws.Range("G17").Value = 47
$ws.Range("H17").Value = -51.063829787234
$ws.Range("I17").Value = 196
$ws.Range("J17").Value = 136
$ws.Range("K17").Value = 44.117647058823
$ws.Range("L17").Value = 13.294797687861
$ws.Range("M17").Value = 157.894736842105
$ws.Range("N17").Value = -26.037735849056
$ws.Range("D18").Value = 16
$ws.Range("E18").Value = -56.25
$ws.Range("F18").Value = 37
$ws.Range("G18").Value = 67
$ws.Range("H18").Value = -44.776119402985
$ws.Range("I18").Value = 181
$ws.Range("J18").Value = 262
$ws.Range("K18").Value = -30.916030534351
$ws.Range("L18").Value = 14.556962025316
$ws.Range("M18").Value = 41.40625
$ws.Range("N18").Value = -82.461240310077
$ws.Range("C19").Value = 54
$ws.Range("D19").Value = 40
$ws.Range("E19").Value = 35
$ws.Range("F19").Value = 203
$ws.Range("G19").Value = 169
$ws.Range("H19").Value = 20.118343195266
$ws.Range("I19").Value = 951
$ws.Range("J19").Value = 769
$ws.Range("K19").Value = 23.667100130039
$ws.Range("L19").Value = 133.660933660934
$ws.Range("M19").Value = 8.685714285714
$ws.Range("N19").Value = -74.199674443841
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 10
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 66.666666666666
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 29
$ws.Range("K20").Value = -10.344827586206
$ws.Range("L20").Value = 36.842105263157
$ws.Range("M20").Value = 188.888888888889
$ws.Range("N20").Value = -83.333333333333
$ws.Range("C21").Value = 75
$ws.Range("D21").Value = 87
$ws.Range("E21").Value = -13.793103448275
$ws.Range("F21").Value = 311
$ws.Range("G21").Value = 351
$ws.Range("H21").Value = -11.396011396011
$ws.Range("I21").Value = 1552
$ws.Range("J21").Value = 1418
$ws.Range("K21").Value = 9.449929478138
$ws.Range("L21").Value = 71.871539313399
$ws.Range("M21").Value = 35.309503051438
$ws.Range("N21").Value = -74.677761461902
$ws.Range("C22").Value = 9
$ws.Range("D22").Value = 6
$ws.Range("E22").Value = 50
$ws.Range("F22").Value = 22
$ws.Range("H22").Value = -4.347826086956
$ws.Range("I22").Value = 95
$ws.Range("J22").Value = 77
$ws.Range("K22").Value = 23.376623376623
$ws.Range("L22").Value = 72.727272727272
$ws.Range("M22").Value = 90
$ws.Range("C24").Value = 79
$ws.Range("D24").Value = 57
$ws.Range("E24").Value = 38.596491228070
$ws.Range("F24").Value = 311
$ws.Range("G24").Value = 256
$ws.Range("H24").Value = 21.484375
$ws.Range("I24").Value = 1435
$ws.Range("J24").Value = 1110
$ws.Range("K24").Value = 29.279279279279
$ws.Range("L24").Value = 82.569974554707
$ws.Range("M24").Value = -19.517666853617
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 25
$ws.Range("E25").Value = -12
$ws.Range("F25").Value = 86
$ws.Range("G25").Value = 64
$ws.Range("H25").Value = 34.375
$ws.Range("I25").Value = 389
$ws.Range("J25").Value = 324
$ws.Range("K25").Value = 20.061728395061
$ws.Range("L25").Value = 23.492063492063
$ws.Range("M25").Value = 62.083333333333
$ws.Range("C26").Value = 1
$ws.Range("D26").Value = 2
$ws.Range("E26").Value = -50
$ws.Range("F26").Value = 3
$ws.Range("G26").Value = 5
$ws.Range("H26").Value = -40
$ws.Range("I26").Value = 9
$ws.Range("J26").Value = 11
$ws.Range("K26").Value = -18.181818181818
$ws.Range("L26").Value = -43.75
$ws.Range("C27").Value = 4
$ws.Range("D27").Value = 4
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 14
$ws.Range("H27").Value = -17.647058823529
$ws.Range("I27").Value = 84
$ws.Range("J27").Value = 75
$ws.Range("K27").Value = 12
$ws.Range("L27").Value = 52.727272727272
$ws.Range("L28").Value = -50
$ws.Range("L29").Value = -50
$ws.Range("L30").Value = -68.75
